$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.616267204284668
$ws.Range("B1").Value = 2.718456983566284
$ws.Range("C1").Value = 3.202774047851562
$ws.Range("D1").Value = 3.477633476257324
$ws.Range("E1").Value = 1.912714719772339
